$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-16 Monday" "2025-06-17 Tuesday"

Replace-Text "133÷8=16, 5" "424÷4=106, 0"
Replace-Text "367÷9=40, 7" "858÷5=171, 3"
Replace-Text "685÷4=171, 1" "527÷7=75, 2"
Replace-Text "370÷4=92, 2" "259÷9=28, 7"
Replace-Text "826÷6=137, 4" "729÷3=243, 0"

Replace-Text "970÷6=161, 4" "312÷9=34, 6"
Replace-Text "488÷2=244, 0" "582÷3=194, 0"
Replace-Text "661÷4=165, 1" "378÷4=94, 2"
Replace-Text "526÷2=263, 0" "430÷6=71, 4"
Replace-Text "555÷6=92, 3" "767÷3=255, 2"

Replace-Text "733÷7=104, 5" "413÷6=68, 5"
Replace-Text "867÷4=216, 3" "416÷4=104, 0"
Replace-Text "914÷3=304, 2" "376÷5=75, 1"
Replace-Text "362÷3=120, 2" "584÷5=116, 4"
Replace-Text "122÷3=40, 2" "688÷5=137, 3"

Replace-Text "594÷3=198, 0" "201÷6=33, 3"
Replace-Text "451÷2=225, 1" "279÷8=34, 7"
Replace-Text "657÷6=109, 3" "142÷8=17, 6"
Replace-Text "123÷7=17, 4" "750÷2=375, 0"
Replace-Text "781÷4=195, 1" "746÷8=93, 2"

Replace-Text "224÷4=56, 0" "539÷8=67, 3"
Replace-Text "971÷3=323, 2" "342÷3=114, 0"
Replace-Text "879÷8=109, 7" "609÷7=87, 0"
Replace-Text "691÷3=230, 1" "399÷6=66, 3"
Replace-Text "930÷9=103, 3" "513÷8=64, 1"
